$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 value from "Teams in Digital" to "Teams in EMU"
$ws.Range("B1").Value = "Teams in EMU"

# Update the active selection to B2 (was A2)
$ws.Range("B2").Select()
